$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'286.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'4.23%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'28.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.08%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.938"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.94%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'2.42%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.244"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'4.52%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'10.56%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9172"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.72%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1567"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'3.30%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06513"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'28.63%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07657"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.74%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.02976"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.47%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08977"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.11%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001604"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.70%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006576"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'2.19%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006082"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.75%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.490"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.60%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.386"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.32%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.244"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.25%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'0.1348"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.56%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.016"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.74%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'12.23%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04458"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.15%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'0.82%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004346"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'12.89%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E28").Value = "'-1.78%"
$ws.Range("E28").Style = "Normal"
$ws.Range("E40").Value = "'0.88%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007069"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.67%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1414"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'20.49%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002049"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-1.07%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'8.71%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005550"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'7.11%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-7.71%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-7.66%"
$ws.Range("E47").Style = "Normal"
